# Updated TPM-derived NATMI ligand-receptor metrics for Jag1-Notch3
# (Ligand avg/total expression for ECs senders and Receptor avg/total
# expression for ECs targets changed; all dependent specificity /
# edge-weight columns recomputed accordingly.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("G2").Value = 11.190867
$ws.Range("H2").Value = 33.572601
$ws.Range("I2").Value = 0.1514016037116739
$ws.Range("J2").Value = 0.1514016037116739
$ws.Range("M2").Value = 7.413580666666667
$ws.Range("N2").Value = 22.240742
$ws.Range("O2").Value = 0.05108888817597561
$ws.Range("P2").Value = 0.05108888817597561
$ws.Range("Q2").Value = 82.964395234438
$ws.Range("R2").Value = 746.679557109942
$ws.Range("S2").Value = 0.00773493960168908
$ws.Range("T2").Value = 0.007734939601689079

# row 3
$ws.Range("G3").Value = 11.190867
$ws.Range("H3").Value = 33.572601
$ws.Range("I3").Value = 0.1514016037116739
$ws.Range("J3").Value = 0.1514016037116739
$ws.Range("O3").Value = 0.0112127179963522
$ws.Range("P3").Value = 0.0112127179963522
$ws.Range("Q3").Value = 18.20858508992
$ws.Range("R3").Value = 163.87726580928
$ws.Range("S3").Value = 0.00169762348661447
$ws.Range("T3").Value = 0.00169762348661447

# row 4
$ws.Range("G4").Value = 11.190867
$ws.Range("H4").Value = 33.572601
$ws.Range("I4").Value = 0.1514016037116739
$ws.Range("J4").Value = 0.1514016037116739
$ws.Range("O4").Value = 0.9376983938276722
$ws.Range("P4").Value = 0.9376983938276722
$ws.Range("Q4").Value = 1522.749524089268
$ws.Range("R4").Value = 13704.74571680341
$ws.Range("S4").Value = 0.1419690406233703
$ws.Range("T4").Value = 0.1419690406233703

# row 5
$ws.Range("I5").Value = 0.2043341870182926
$ws.Range("J5").Value = 0.2043341870182926
$ws.Range("M5").Value = 7.413580666666667
$ws.Range("N5").Value = 22.240742
$ws.Range("O5").Value = 0.05108888817597561
$ws.Range("P5").Value = 0.05108888817597561
$ws.Range("Q5").Value = 111.9701630372233
$ws.Range("R5").Value = 1007.73146733501
$ws.Range("S5").Value = 0.01043920643110644
$ws.Range("T5").Value = 0.01043920643110644

# row 6
$ws.Range("I6").Value = 0.2043341870182926
$ws.Range("J6").Value = 0.2043341870182926
$ws.Range("O6").Value = 0.0112127179963522
$ws.Range("P6").Value = 0.0112127179963522
$ws.Range("S6").Value = 0.002291141616050006
$ws.Range("T6").Value = 0.002291141616050006

# row 7
$ws.Range("I7").Value = 0.2043341870182926
$ws.Range("J7").Value = 0.2043341870182926
$ws.Range("O7").Value = 0.9376983938276722
$ws.Range("P7").Value = 0.9376983938276722
$ws.Range("S7").Value = 0.1916038389711362
$ws.Range("T7").Value = 0.1916038389711361

# row 8
$ws.Range("G8").Value = 47.62086333333333
$ws.Range("I8").Value = 0.6442642092700336
$ws.Range("J8").Value = 0.6442642092700336
$ws.Range("M8").Value = 7.413580666666667
$ws.Range("N8").Value = 22.240742
$ws.Range("O8").Value = 0.05108888817597561
$ws.Range("P8").Value = 0.05108888817597561
$ws.Range("Q8").Value = 353.0411117379755
$ws.Range("R8").Value = 3177.37000564178
$ws.Range("S8").Value = 0.0329147421431801
$ws.Range("T8").Value = 0.0329147421431801

# row 9
$ws.Range("G9").Value = 47.62086333333333
$ws.Range("I9").Value = 0.6442642092700336
$ws.Range("J9").Value = 0.6442642092700336
$ws.Range("O9").Value = 0.0112127179963522
$ws.Range("P9").Value = 0.0112127179963522
$ws.Range("Q9").Value = 77.48358925724445
$ws.Range("R9").Value = 697.3523033151999
$ws.Range("S9").Value = 0.007223952893687726
$ws.Range("T9").Value = 0.007223952893687726

# row 10
$ws.Range("G10").Value = 47.62086333333333
$ws.Range("I10").Value = 0.6442642092700336
$ws.Range("J10").Value = 0.6442642092700336
$ws.Range("O10").Value = 0.9376983938276722
$ws.Range("P10").Value = 0.9376983938276722
$ws.Range("Q10").Value = 6479.805986216564
$ws.Range("R10").Value = 58318.25387594908
$ws.Range("S10").Value = 0.6041255142331657
$ws.Range("T10").Value = 0.6041255142331657
